# Auto-generated Excel COM-interop edit script
# Applies the per-cell value changes described by the commit diff
# ("Code updated 23-04-27 11:30:55") to Season_Trophies/89.xlsx (Sheet1).
#
# Columns A, B and E hold numeric-looking values that must remain stored
# as TEXT (the source file uses inlineStr for every cell). Assigning a
# plain numeric-looking string via .Value/.Value2 gets auto-coerced to a
# real number, so we use the classic Excel "quote-prefix" trick (a leading
# apostrophe via .Value2) to force text storage, matching the original type.
# Column C (Name) is plain text and never looks numeric, so it is assigned
# directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = '''59123'
$ws.Range("E2").Value2 = '''2534'
$ws.Range("A3").Value2 = '''61117'
$ws.Range("E3").Value2 = '''2515'
$ws.Range("A5").Value2 = '''40179'
$ws.Range("E5").Value2 = '''3067'
$ws.Range("A7").Value2 = '''61041'
$ws.Range("A8").Value2 = '''7211'
$ws.Range("E8").Value2 = '''4206'
$ws.Range("A9").Value2 = '''8445'
$ws.Range("B9").Value2 = '''53060417'
$ws.Range("C9").Value = '㊥老纳信耶稣'
$ws.Range("E9").Value2 = '''4157'
$ws.Range("A10").Value2 = '''10533'
$ws.Range("B10").Value2 = '''38809086'
$ws.Range("C10").Value = 'Kouenᶻᵍˣ'
$ws.Range("E10").Value2 = '''4072'
$ws.Range("A11").Value2 = '''12353'
$ws.Range("B11").Value2 = '''53520939'
$ws.Range("C11").Value = '㊥虎哥tiger'
$ws.Range("E11").Value2 = '''4007'
$ws.Range("A12").Value2 = '''14308'
$ws.Range("B12").Value2 = '''4756174'
$ws.Range("C12").Value = '純希です'
$ws.Range("E12").Value2 = '''3985'
$ws.Range("A13").Value2 = '''36024'
$ws.Range("E13").Value2 = '''3254'
$ws.Range("A14").Value2 = '''57027'
$ws.Range("A15").Value2 = '''62927'
$ws.Range("A17").Value2 = '''12287'
$ws.Range("B17").Value2 = '''46289694'
$ws.Range("C17").Value = '㊥Vincent'
$ws.Range("E17").Value2 = '''4009'
$ws.Range("A18").Value2 = '''12462'
$ws.Range("B18").Value2 = '''8057001'
$ws.Range("C18").Value = '㊥兵者诡道也'
$ws.Range("E18").Value2 = '''4004'
$ws.Range("A19").Value2 = '''13015'
$ws.Range("B19").Value2 = '''31134300'
$ws.Range("C19").Value = 'McMaX'
$ws.Range("E19").Value2 = '''3996'
$ws.Range("A20").Value2 = '''13107'
$ws.Range("B20").Value2 = '''54698813'
$ws.Range("C20").Value = '閃亮唐老鴨'
$ws.Range("E20").Value2 = '''3995'
$ws.Range("A21").Value2 = '''13764'
$ws.Range("B21").Value2 = '''31495601'
$ws.Range("C21").Value = '陈晓军'
$ws.Range("E21").Value2 = '''3991'
$ws.Range("A22").Value2 = '''17853'
$ws.Range("B22").Value2 = '''54085771'
$ws.Range("C22").Value = '㊥Matthieu'
$ws.Range("E22").Value2 = '''3867'
$ws.Range("A23").Value2 = '''20581'
$ws.Range("B23").Value2 = '''55769051'
$ws.Range("C23").Value = '㊥叮叮当.'
$ws.Range("E23").Value2 = '''3777'
$ws.Range("A24").Value2 = '''22559'
$ws.Range("B24").Value2 = '''56732705'
$ws.Range("C24").Value = '时间温柔皆遗憾'
$ws.Range("E24").Value2 = '''3715'
$ws.Range("A25").Value2 = '''23187'
$ws.Range("B25").Value2 = '''56585361'
$ws.Range("C25").Value = '"㊥ go策划我要ali"'
$ws.Range("E25").Value2 = '''3696'
$ws.Range("A26").Value2 = '''26454'
$ws.Range("E26").Value2 = '''3592'
$ws.Range("A27").Value2 = '''27489'
$ws.Range("E27").Value2 = '''3557'
$ws.Range("A28").Value2 = '''29478'
$ws.Range("E28").Value2 = '''3491'
$ws.Range("A29").Value2 = '''35146'
$ws.Range("E29").Value2 = '''3294'
$ws.Range("A30").Value2 = '''44250'
$ws.Range("E30").Value2 = '''2889'
$ws.Range("A31").Value2 = '''1153'
$ws.Range("E31").Value2 = '''4591'
$ws.Range("A32").Value2 = '''7136'
$ws.Range("E32").Value2 = '''4209'
$ws.Range("A33").Value2 = '''7688'
$ws.Range("B33").Value2 = '''11582001'
$ws.Range("C33").Value = 'iMinatoX4'
$ws.Range("E33").Value2 = '''4189'
$ws.Range("A34").Value2 = '''10979'
$ws.Range("B34").Value2 = '''7852598'
$ws.Range("C34").Value = 'seiji'
$ws.Range("E34").Value2 = '''4053'
$ws.Range("A35").Value2 = '''11108'
$ws.Range("B35").Value2 = '''35114520'
$ws.Range("C35").Value = '13lur¹³'
$ws.Range("E35").Value2 = '''4049'
$ws.Range("A36").Value2 = '''13343'
$ws.Range("B36").Value2 = '''26280580'
$ws.Range("C36").Value = '꧁SSS.TIGRESS꧂ᶻᵍˣ'
$ws.Range("E36").Value2 = '''3994'
$ws.Range("A37").Value2 = '''13939'
$ws.Range("B37").Value2 = '''55317038'
$ws.Range("C37").Value = 'necman12345'
$ws.Range("E37").Value2 = '''3989'
$ws.Range("A38").Value2 = '''16148'
$ws.Range("B38").Value2 = '''38995116'
$ws.Range("C38").Value = '"Ramesh Pavai Nam"'
$ws.Range("E38").Value2 = '''3928'
$ws.Range("A39").Value2 = '''18062'
$ws.Range("E39").Value2 = '''3859'
$ws.Range("A40").Value2 = '''18981'
$ws.Range("B40").Value2 = '''45967307'
$ws.Range("C40").Value = 'Ricky'
$ws.Range("E40").Value2 = '''3829'
$ws.Range("A41").Value2 = '''20933'
$ws.Range("B41").Value2 = '''6809364'
$ws.Range("C41").Value = '"Scorp IP"'
$ws.Range("E41").Value2 = '''3767'
$ws.Range("A42").Value2 = '''25789'
$ws.Range("B42").Value2 = '''56379103'
$ws.Range("C42").Value = 'Globalking'
$ws.Range("E42").Value2 = '''3613'
$ws.Range("A43").Value2 = '''26185'
$ws.Range("B43").Value2 = '''47459684'
$ws.Range("C43").Value = '㊥阿闹切克闹'
$ws.Range("E43").Value2 = '''3600'
$ws.Range("A44").Value2 = '''27138'
$ws.Range("E44").Value2 = '''3568'
$ws.Range("A45").Value2 = '''32715'
$ws.Range("B45").Value2 = '''50837459'
$ws.Range("C45").Value = 'NINE日'
$ws.Range("E45").Value2 = '''3382'
$ws.Range("A46").Value2 = '''35804'
$ws.Range("B46").Value2 = '''58203298'
$ws.Range("C46").Value = '权旨qua'
$ws.Range("E46").Value2 = '''3264'
$ws.Range("A47").Value2 = '''35901'
$ws.Range("B47").Value2 = '''52997727'
$ws.Range("C47").Value = 'larios'
$ws.Range("E47").Value2 = '''3260'
$ws.Range("A48").Value2 = '''40011'
$ws.Range("E48").Value2 = '''3076'
$ws.Range("A49").Value2 = '''40106'
$ws.Range("E49").Value2 = '''3071'
$ws.Range("A50").Value2 = '''52001'
$ws.Range("A51").Value2 = '''52783'
$ws.Range("E51").Value2 = '''2633'
$ws.Range("A52").Value2 = '''55578'
$ws.Range("E52").Value2 = '''2581'
$ws.Range("A53").Value2 = '''61062'
$ws.Range("B53").Value2 = '''1550355'
$ws.Range("C53").Value = '"皓茵 世界"'
$ws.Range("E53").Value2 = '''2516'
$ws.Range("A54").Value2 = '''63069'
$ws.Range("B54").Value2 = '''20737010'
$ws.Range("C54").Value = '混着玩...'
$ws.Range("E54").Value2 = '''2500'
$ws.Range("A56").Value2 = '''54126'
$ws.Range("E56").Value2 = '''2606'
$ws.Range("A57").Value2 = '''58875'
$ws.Range("E57").Value2 = '''2537'
$ws.Range("A60").Value2 = '''34816'
$ws.Range("E60").Value2 = '''3306'
$ws.Range("A61").Value2 = '''38606'
$ws.Range("E61").Value2 = '''3140'
$ws.Range("A62").Value2 = '''56194'
$ws.Range("E62").Value2 = '''2571'
$ws.Range("A63").Value2 = '''63158'
$ws.Range("A64").Value2 = '''63177'
$ws.Range("B64").Value2 = '''9718882'
$ws.Range("C64").Value = '小霸王2021'
$ws.Range("E64").Value2 = '''2499'
$ws.Range("A65").Value2 = '''88106'
$ws.Range("B65").Value2 = '''57219176'
$ws.Range("C65").Value = '青莲道人'
$ws.Range("E65").Value2 = '''1522'
$ws.Range("A66").Value2 = '''101772'
$ws.Range("B66").Value2 = '''49000199'
$ws.Range("C66").Value = 'SlipperyForester5672'
$ws.Range("E66").Value2 = '''1284'
$ws.Range("B67").Value2 = '''56700848'
$ws.Range("C67").Value = '工口漫画老师'
$ws.Range("B68").Value2 = '''38994054'
$ws.Range("C68").Value = 'chengnan'
$ws.Range("B69").Value2 = '''3391765'
$ws.Range("C69").Value = '马er'
$ws.Range("B70").Value2 = '''55810157'
$ws.Range("C70").Value = 'Beard'
$ws.Range("B71").Value2 = '''57556179'
$ws.Range("C71").Value = '特战新生代英雄'
$ws.Range("B72").Value2 = '''1222440'
$ws.Range("C72").Value = '"Sneaky Ninja Panda"'
$ws.Range("B73").Value2 = '''58340439'
$ws.Range("C73").Value = '70qilin'
$ws.Range("B74").Value2 = '''15436348'
$ws.Range("C74").Value = 'Lucas'
$ws.Range("B75").Value2 = '''20372140'
$ws.Range("C75").Value = '人山即是仙'
$ws.Range("B76").Value2 = '''58615925'
$ws.Range("C76").Value = '齐天的大圣'
$ws.Range("B77").Value2 = '''58641574'
$ws.Range("C77").Value = 'Player-58641574鱼'
$ws.Range("B78").Value2 = '''58743790'
$ws.Range("C78").Value = 'Ma'
$ws.Range("B79").Value2 = '''54941706'
$ws.Range("C79").Value = 'AlexMenjivar20'
$ws.Range("A80").Value2 = '''46207'
$ws.Range("E80").Value2 = '''2821'
$ws.Range("A83").Value2 = '''140844'
